$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.489.23'
$ws.Range("E2").Value = '  -3.29%  '

$ws.Range("D3").Value = '2.984.55'
$ws.Range("E3").Value = '  -5.10%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''496.60'
$ws.Range("E5").Value = '  -5.79%  '

$ws.Range("D6").Value = '''135.11'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '2.981.93'
$ws.Range("E8").Value = '  -5.09%  '

$ws.Range("E9").Value = '  -3.63%  '

$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("E11").Value = '  -3.55%  '

$ws.Range("E12").Value = '  -6.91%  '

$ws.Range("E13").Value = '  -0.48%  '

$ws.Range("D14").Value = '3.495.31'
$ws.Range("E14").Value = '  -4.99%  '

$ws.Range("D15").Value = '''25.05'
$ws.Range("E15").Value = '  -2.06%  '

$ws.Range("D16").Value = '56.403.38'
$ws.Range("E16").Value = '  -3.21%  '

$ws.Range("D17").Value = '2.983.49'
$ws.Range("E17").Value = '  -4.77%  '

$ws.Range("E18").Value = '  -4.47%  '

$ws.Range("D19").Value = '''5.83'
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").Value = '''12.38'
$ws.Range("E20").Value = '  -5.15%  '

$ws.Range("E21").Value = '  -2.39%  '

$ws.Range("D22").Value = '''325.50'
$ws.Range("E22").Value = '  -5.66%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("E24").Value = '  -7.97%  '

$ws.Range("D25").Value = '''61.38'
$ws.Range("E25").Value = '  -9.58%  '

$ws.Range("D26").Value = '''0.997'
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("E27").Value = '  -4.29%  '

$ws.Range("D28").Value = '0.0₃0899'
$ws.Range("E28").Value = '  -6.56%  '

$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("E30").Value = '  -4.49%  '

$ws.Range("D31").Value = '''6.66'
$ws.Range("E31").Value = '  -3.00%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '''1.17'
$ws.Range("E32").Value = '  -4.20%  '

$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '''1.74'
$ws.Range("E33").Value = '  -6.76%  '

$ws.Range("D34").Value = '''20.19'
$ws.Range("E34").Value = '  -5.79%  '

$ws.Range("D35").Value = '''154.94'
$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("D36").Value = '''4.47'
$ws.Range("E36").Value = '  -7.20%  '

$ws.Range("E37").Value = '  -6.39%  '

$ws.Range("E38").Value = '  -10.45%  '

$ws.Range("D39").Value = '''0.0685'
$ws.Range("E39").Value = '  -0.29%  '

$ws.Range("D40").Value = '''23.10'
$ws.Range("E40").Value = '  -4.85%  '

$ws.Range("D41").Value = '3.018.04'
$ws.Range("E41").Value = '  -4.85%  '

$ws.Range("D42").Value = '''36.58'
$ws.Range("E42").Value = '  -9.61%  '

$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  +0.06%  '

$ws.Range("E44").Value = '  -7.72%  '

$ws.Range("D45").Value = '2.234.17'
$ws.Range("E45").Value = '  -1.72%  '

$ws.Range("D46").Value = '''0.995'
$ws.Range("E46").Value = '  -8.65%  '

$ws.Range("E47").Value = '  -2.92%  '

$ws.Range("E48").Value = '  -8.64%  '

$ws.Range("D49").Value = '''1.95'
$ws.Range("E49").Value = '  +3.37%  '

$ws.Range("E50").Value = '  +1.41%  '

$ws.Range("D51").Value = '''5.79'
$ws.Range("E51").Value = '  -7.01%  '
